# Restored from revision of admin on 04/24/2020 09:01:59 AM.TEST Author: admin. Type: SAVE.
# Change: cell C10 ("From" value of rule R30) changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10").Value = 1
